$wb = $excel.ActiveWorkbook

# "展览" (Exhibition) sheet holds the per-category rows.
$sheetExhibition = $wb.Worksheets.Item("展览")
# "全部类型" (All types) sheet aggregates rows from every category sheet.
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# 苏州·世纪幻想动漫游戏展2.0 -- 想去人数 (want-to-go count) 1440 -> 1442
$sheetExhibition.Range("F4").Value  = 1442
# 【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会 -- 10812 -> 10813
$sheetExhibition.Range("F7").Value  = 10813
# 苏州·第十七届 I COME ACG 动漫品牌博览会 -- 12590 -> 12591
$sheetExhibition.Range("F14").Value = 12591

# Same three events mirrored on the aggregated "全部类型" sheet, one row lower.
$sheetAllTypes.Range("F5").Value  = 1442
$sheetAllTypes.Range("F8").Value  = 10813
$sheetAllTypes.Range("F15").Value = 12591
